# Update header row (row 1) on the active sheet to the new reference-level
# column labels, and relabel the "R" row marker.
#
# Old headers:  typ | min | SD_nedre | SD_D | D_M | M_G | G_SG | SG_øvre | max
# New headers:  typ | pess | X0 | X20 | X40 | X60 | X80 | X100 | opt
#
# Column A keeps "typ" in row 1 and the data row (row 2) keeps its "R" marker;
# only the shared-string pool order changes as a side effect of the relabel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "pess"
$ws.Range("C1").Value = "X0"
$ws.Range("D1").Value = "X20"
$ws.Range("E1").Value = "X40"
$ws.Range("F1").Value = "X60"
$ws.Range("G1").Value = "X80"
$ws.Range("H1").Value = "X100"
$ws.Range("I1").Value = "opt"

$ws.Range("A2").Value = "R"
